$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Team1 row (row 2): name casing + GroupNo
$ws.Range("B2").Value = "Team1"
$ws.Range("D2").Value = 1

# Update row 3: registration date + GroupNo, then clear score/game stats
$ws.Range("C3").Value = "15/12"
$ws.Range("D3").Value = 2
$ws.Range("E3:I3").ClearContents()

# Remove the now-obsolete team3/team4/team5 rows (old rows 4-6)
$ws.Range("A4:I6").EntireRow.Delete()
